# Orders.xlsx update:
#  - Row 9 becomes Rhonbrau Klosterbier / 3
#  - Row 10 becomes "Not a real beverage" / 4 (typo "lol" dropped from the product name)
#  - Old rows 11 and 12 (which used to hold these two products) are removed
#  - Selection moves to A12 (now an empty cell below the shrunk table)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Fix the typo'd product name everywhere it is used in the sheet.
$ws.Range("A12").Value2 = "Not a real beverage"

# Move the last two order rows up onto rows 9/10, overwriting the
# superseded "Laughing Lumberjack Lager" / "Outback Lager" duplicate rows.
$ws.Range("A9").Value2 = $ws.Range("A11").Value2
$ws.Range("B9").Value2 = $ws.Range("B11").Value2
$ws.Range("A10").Value2 = $ws.Range("A12").Value2
$ws.Range("B10").Value2 = $ws.Range("B12").Value2

# Remove the now-duplicated trailing rows 11:12 so the table shrinks back
# down to A1:D10.
$ws.Range("A11:D12").Delete()

$ws.Range("A12").Select()
